$wb = $excel.ActiveWorkbook

# --- Remove the empty "Planilha2" worksheet (it was added by mistake) ---
$excel.DisplayAlerts = $false
[void]$wb.Worksheets.Item("Planilha2").Delete()
$excel.DisplayAlerts = $true

$ws = $wb.Worksheets.Item("Planilha1")

# --- Content fixes on the weekly schedule table ---
# B8 was "Estudos", B9 was "Jantar" -> both become "Igreja"
$ws.Range("B8").Value = "Igreja"
$ws.Range("B9").Value = "Igreja"

# --- Remove now-unneeded direct formatting on F6:G10 (Qui/Sex columns) ---
$ws.Range("F6:G10").ClearFormats()

# --- Conditional formatting had a duplicated "last week" time-period rule;
#     drop the extra one and keep the other as the sole, top-priority rule ---
$fcs = $ws.Range("A3:H10").FormatConditions
$fcs.Item(1).Delete()
$fcs.Item(1).SetFirstPriority()

# --- Restore cursor/selection to where the author left off ---
$ws.Range("D16").Select()
